# Daily attendance processing - 2025-12-25 04:25:10
# Reverses the order of comma-separated entries in the "Recorded By" (column G)
# values on the "Session Analysis Results" sheet, wherever more than one
# entry is present. Single-entry cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text.Split(",")
    $count = $parts.Count

    if ($count -gt 1) {
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $reversed = $trimmed[($count - 1)..0]
        $joined = [string]::Join(", ", $reversed)

        if ($joined -ne $text) {
            $cell.Value = $joined
        }
    }
}
